$wb = $excel.ActiveWorkbook

# --- Add the new "positionData" worksheet as the last (3rd) tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "positionData"

# --- Header row (row 3) ---
$ws.Range("B3").Value = "D"
$ws.Range("C3").Value = "M"
$ws.Range("D3").Value = "M’"
$ws.Range("E3").Value = "F"
$ws.Range("F3").Value = "sine coefficient"
$ws.Range("G3").Value = "cosine coefficient"

$ws.Range("B3:E3").HorizontalAlignment = -4108
$ws.Range("B3:E3").Font.Bold = $true
$ws.Range("F3:G3").Font.Bold = $true

# --- Data rows (4-35) ---
$data = @(
    @(0,0,1,0,6288774,-20905355),
    @(2,0,-1,0,1274027,-1699111),
    @(2,0,0,0,$null,$null),
    @(0,0,2,0,$null,$null),
    @(0,1,0,0,$null,$null),
    @(0,0,0,2,$null,$null),
    @(2,0,-2,0,$null,$null),
    @(2,-1,-1,0,$null,$null),
    @(2,0,1,0,$null,$null),
    @(2,-1,0,0,$null,$null),
    @(0,1,-1,0,$null,$null),
    @(1,0,0,0,$null,$null),
    @(0,1,1,0,$null,$null),
    @(2,0,0,-2,$null,$null),
    @(0,0,1,2,$null,$null),
    @(0,0,1,-2,$null,$null),
    @(4,0,-1,0,$null,$null),
    @(0,0,3,0,$null,$null),
    @(4,0,-2,0,$null,$null),
    @(2,1,-1,0,$null,$null),
    @(2,1,0,0,$null,$null),
    @(1,0,-1,0,$null,$null),
    @(1,1,0,0,$null,$null),
    @(2,-1,1,0,$null,$null),
    @(2,0,2,0,$null,$null),
    @(4,0,0,0,$null,$null),
    @(2,0,-3,0,$null,$null),
    @(0,1,-2,0,$null,$null),
    @(2,0,-1,2,$null,$null),
    @(2,-1,-2,0,$null,$null),
    @(1,0,1,0,$null,$null),
    @(2,-2,0,0,$null,$null)
)

$row = 4
foreach ($r in $data) {
    $ws.Cells.Item($row, 2).Value = $r[0]
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    if ($r[4] -ne $null) {
        $ws.Cells.Item($row, 6).Value = $r[4]
    }
    if ($r[5] -ne $null) {
        $ws.Cells.Item($row, 7).Value = $r[5]
    }
    $row = $row + 1
}

# --- Column widths ---
$ws.Columns.Item(6).ColumnWidth = 13.65
$ws.Columns.Item(7).ColumnWidth = 15.61

Write-Output "done"
